$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 40-41, shifting the existing rows 40-52 down to 42-54.
$ws.Rows("40:41").Insert()

# Row 40: new weekly price entry - "Primera" quality
$ws.Range("A40").Value = 1
$ws.Range("B40").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C40").Value = "Arica y Parinacota"
$ws.Range("D40").Value = 44438
$ws.Range("E40").Value = 15
$ws.Range("F40").Value = 100112036
$ws.Range("G40").Value = "Caigua"
$ws.Range("H40").Value = "Sin especificar"
$ws.Range("I40").Value = "Primera"
$ws.Range("J40").Value = 120
$ws.Range("K40").Value = 6000
$ws.Range("L40").Value = 7000
$ws.Range("M40").Value = 6500
$ws.Range("N40").Value = "$/caja 20 kilos"
$ws.Range("O40").Value = "Región de Arica y Parinacota"
$ws.Range("P40").Value = 325
$ws.Range("Q40").Value = 20
$ws.Range("R40").Value = "Hortaliza"

# Row 41: new weekly price entry - "Segunda" quality
$ws.Range("A41").Value = 1
$ws.Range("B41").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C41").Value = "Arica y Parinacota"
$ws.Range("D41").Value = 44438
$ws.Range("E41").Value = 15
$ws.Range("F41").Value = 100112036
$ws.Range("G41").Value = "Caigua"
$ws.Range("H41").Value = "Sin especificar"
$ws.Range("I41").Value = "Segunda"
$ws.Range("J41").Value = 120
$ws.Range("K41").Value = 5000
$ws.Range("L41").Value = 6000
$ws.Range("M41").Value = 5500
$ws.Range("N41").Value = "$/caja 20 kilos"
$ws.Range("O41").Value = "Región de Arica y Parinacota"
$ws.Range("P41").Value = 275
$ws.Range("Q41").Value = 20
$ws.Range("R41").Value = "Hortaliza"
